$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Agency name in column A (rows 2-42) from
# "Judicial Branch - CSSD" to "Judicial Branch - Court Support Services Division (JB-CSSD)"
for ($r = 2; $r -le 42; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Text -eq "Judicial Branch - CSSD") {
        $cell.Value = "Judicial Branch - Court Support Services Division (JB-CSSD)"
    }
}

# Move the active selection to C8
$ws.Range("C8").Select() | Out-Null
